$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the three address hyperlinks with new email addresses ---
# Remove the existing (stale) hyperlinks on the sheet first so old
# mailto: targets for aaa/bbb/ccc@sample.xyz don't linger.
$ws.Range("A2").Hyperlinks.Delete()

# Write the new e-mail text in row order (A2, A3, A4) so the shared-string
# table is built up in that same order.
$ws.Range("A2").Value = "gianluca.marchitelli14@gmail.com"
$ws.Range("A3").Value = "hello@johnnybrr.site"
$ws.Range("A4").Value = "gianluca_marchitelli@yahoo.it"

# Re-create the hyperlinks (this also applies the built-in "Hyperlink"
# cell style - underline + theme font color - to each cell). Add them in
# the order A2, A4, A3 so the generated relationship ids land as
# rId1->A2, rId2->A4, rId3->A3.
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:gianluca.marchitelli14@gmail.com", [Type]::Missing, [Type]::Missing, "gianluca.marchitelli14@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:gianluca_marchitelli@yahoo.it", [Type]::Missing, [Type]::Missing, "gianluca_marchitelli@yahoo.it") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:hello@johnnybrr.site", [Type]::Missing, [Type]::Missing, "hello@johnnybrr.site") | Out-Null

# --- Prep the sheet for the new per-message send loop/log ---
# Reserve a styled (but empty) row far below the existing data, matching
# the "Hyperlink" look, so the loop has a formatted cell to grow into.
$ws.Range("A252").Style = "Hyperlink"

# Reflect the in-progress multi-row selection the author left active.
$ws.Range("A3:A4").Select()

# Restore the page setup for printing the address list.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
